# Apply cryptos list update (prices + 1h volume %) per commit:
# "Updated cryptos list on Tue Aug 22 03:11:45 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.137.53'
$ws.Range('E2').Value = '  -0.69%  '
$ws.Range('D3').Value = '1.668.88'
$ws.Range('E4').Value = '  -0.54%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '210.58'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -3.37%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5252'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -2.77%  '
$ws.Range('E7').Value = '  -0.51%  '
$ws.Range('E8').Value = '  -3.70%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06305'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -2.31%  '
$ws.Range('E10').Value = '  -2.36%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07555'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -1.42%  '
$ws.Range('D12').Value = '1.661.62'
$ws.Range('E12').Value = '  -5.23%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.447'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -2.11%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.5564'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -3.97%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '66.88'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.03%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.000007934'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -5.61%  '
$ws.Range('D17').Value = '26.164.98'
$ws.Range('E17').Value = '  -0.82%  '
$ws.Range('E18').Value = '  -0.54%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '4.753'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -3.29%  '
$ws.Range('B20').Value = 'BitcoinCash'
$ws.Range('C20').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '186.74'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.92%  '
$ws.Range('B21').Value = 'Avalanche'
$ws.Range('C21').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '10.37'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -4.69%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.175'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.37%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.003'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.56%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '149.78'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.43%  '
$ws.Range('E25').Value = '  -2.81%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.507'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -4.62%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '15.97'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.89%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.06247'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.78%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.353'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -1.80%  '
$ws.Range('E30').Value = '  -3.50%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.522'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -2.38%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.420'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -4.56%  '
$ws.Range('E33').Value = '  -2.64%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.9983'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -3.27%  '
$ws.Range('E35').Value = '  -1.91%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.411'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.22%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.732'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -1.45%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '6.144'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.29%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01616'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -2.41%  '
$ws.Range('D40').Value = '1.102.00'
$ws.Range('E40').Value = '  -0.93%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.8735'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -1.27%  '
$ws.Range('E42').Value = '  -0.97%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '100.25'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.96%  '
$ws.Range('D44').Value = '1.822.82'
$ws.Range('E44').Value = '  -1.14%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.00000000107'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -4.40%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '55.49'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -3.93%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.006'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.18%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '8.074'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.42%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.05237'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.90%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.4248'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -1.20%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '5.980'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.50%  '
